$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A19").Value = "EM IMI (EM)"
$ws.Range("B19").Value = 664220
$ws.Range("C19").Value = "EM"
$ws.Range("D19").Value = "IE00BKM4GZ66"
$ws.Range("E19").Value = "ishares"
$ws.Range("G19").Value = 0.18

$ws.Range("G20").Select()
